$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Insert()
$ws.Range("B1").Value = "Country (Population)"
$ws.Range("B2").Select()
